$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("actives_misc")

# Row 30: G30 standalone formula, H30:L30 shared formula (row9 * row19)
$ws.Range("G30").Formula = "=G9*G19"
$ws.Range("H30:L30").Formula = "=H9*H19"

# Rows 31-39: shared formula block (row[n-21] * row[n-11]) across G:L
$ws.Range("G31:L39").Formula = "=G10*G20"

# Row 41: grand total
$ws.Range("G41").Formula = "=SUM(G30:L39)"

# Clear any auto-applied formatting so these new cells keep the default style
$ws.Range("G30:L39").Style = "Normal"
$ws.Range("G41").Style = "Normal"

# Column G width: Excel auto-fit after the new wider salary*headcount totals
# (ColumnWidth chosen so the stored OOXML width lands on the same value a
# real AutoFit on the 9-digit grand-total would produce: 10)
$ws.Columns.Item(7).ColumnWidth = 9.17

# Update the view: scroll down to the new data and select the grand total cell
$ws.Range("G41").Select()
